$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateCell($addr, $value) {
    $cell = $ws.Range($addr)
    $existing = $cell.Value()
    if ($existing -eq $null) {
        # Brand-new cell: copy number format (date style) from a known formatted date cell first.
        $ws.Range("C2").Copy()
        $cell.PasteSpecial(-4122)
    }
    $cell.Value = $value
}

# Row 2 (MVI001R019): add H2, I2
Set-DateCell "H2" 44195
Set-DateCell "I2" 44237

# Row 3 (MVI002R004): add F3, G3, H3, I3
Set-DateCell "F3" 43924
Set-DateCell "G3" 43950
Set-DateCell "H3" 43982
Set-DateCell "I3" 44284

# Row 4 (MVI003R140): update F4 value
$ws.Range("F4").Value = 44263

# Row 5 (MVI004R201): add G5
Set-DateCell "G5" 44271

# Row 6 (MVI005R107): add F6
Set-DateCell "F6" 44249

# Row 7 (MVI006R296): add F7
Set-DateCell "F7" 43942

# Row 8 (MVI007R765): update E8, clear F8 (keep date style, blank value)
$ws.Range("E8").Value = 44312
$ws.Range("F8").ClearContents()

# Row 9 (MVI008R021): add I9
Set-DateCell "I9" 44300

# Row 10 (MVI009R908): add F10, G10
Set-DateCell "F10" 44252
Set-DateCell "G10" 44278

# Update the active selection to H10 as in the target sheet view.
$ws.Range("H10").Select()
